$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Germany Verbandsliga")

# The rows for id=3 (row 5) and id=4 (row 6) had their match data swapped
# (home team identity, away team, and all odds columns B,G,H-AC), while
# the rank column (A) and Div/Div Original Name/Date (C,D,E) stayed put.

# --- Swap HomeTeam (column F) ---
$ws.Range("F5").Value = "SV Altldersdorf"
$ws.Range("F6").Value = "FSV Saxonia Tangermunde"

# --- Swap id (column B) ---
$ws.Range("B5").Value = 6781300
$ws.Range("B6").Value = 6781301

# --- Swap AwayTeam (column G) ---
$ws.Range("G5").Value = "SV Frankonia Wernsdorf"
$ws.Range("G6").Value = "MSC Preussen 1899"

# --- Swap FTHG / FTAG (columns H, I) ---
$ws.Range("H5").Value = 8
$ws.Range("I5").Value = 2
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 0

# --- Swap odds columns K through AC ---
$ws.Range("K5").Value = 2.2
$ws.Range("L5").Value = 3.5
$ws.Range("M5").Value = 2.7
$ws.Range("N5").Value = 1.727
$ws.Range("O5").Value = 4
$ws.Range("P5").Value = 3.5
$ws.Range("Q5").Value = -0.5
$ws.Range("R5").Value = 1.775
$ws.Range("S5").Value = 2.025
$ws.Range("T5").Value = 3.25
$ws.Range("U5").Value = 1.925
$ws.Range("V5").Value = 1.875
$ws.Range("W5").Value = 0.7270000000000001
$ws.Range("X5").Value = -1
$ws.Range("Y5").Value = -1
$ws.Range("Z5").Value = 0.7749999999999999
$ws.Range("AA5").Value = -1
$ws.Range("AB5").Value = 0.925
$ws.Range("AC5").Value = -1

$ws.Range("K6").Value = 1.65
$ws.Range("L6").Value = 4
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 1.45
$ws.Range("O6").Value = 4.5
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = -1.25
$ws.Range("R6").Value = 1.925
$ws.Range("S6").Value = 1.875
$ws.Range("T6").Value = 3.5
$ws.Range("U6").Value = 1.825
$ws.Range("V6").Value = 1.975
$ws.Range("W6").Value = 0.45
$ws.Range("X6").Value = -1
$ws.Range("Y6").Value = -1
$ws.Range("Z6").Value = 0.925
$ws.Range("AA6").Value = -1
$ws.Range("AB6").Value = 0.825
$ws.Range("AC6").Value = -1
